$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previous periods (rows 16-22) are reordered descending (2001 .. 1907),
# carrying their matching "Salario Basico" (column F) along, and the
# outstanding "Valor Mora" (column G) is cleared to 0 for every period
# since those old account statements are no longer due.

$ws.Range("E16").Value = "2001"
$ws.Range("F16").Value = 36120
$ws.Range("G16").Value = 0

$ws.Range("E17").Value = "1912"
$ws.Range("F17").Value = 51600
$ws.Range("G17").Value = 0

$ws.Range("E18").Value = "1911"
$ws.Range("F18").Value = 51600
$ws.Range("G18").Value = 0

$ws.Range("E19").Value = "1910"
$ws.Range("F19").Value = 51600
$ws.Range("G19").Value = 0

$ws.Range("E20").Value = "1909"
$ws.Range("F20").Value = 51600
$ws.Range("G20").Value = 0

$ws.Range("E21").Value = "1908"
$ws.Range("F21").Value = 51600
$ws.Range("G21").Value = 0

$ws.Range("E22").Value = "1907"
$ws.Range("F22").Value = 51600
$ws.Range("G22").Value = 0
